# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 15763
$ws1.Range("F6").Value  = 419
$ws1.Range("F9").Value  = 15460
$ws1.Range("F11").Value = 9065
$ws1.Range("F12").Value = 389
$ws1.Range("F15").Value = 102
$ws1.Range("F18").Value = 204
$ws1.Range("F35").Value = 259

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 70

# Sheet "全部类型" (all types, combined)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 15763
$ws4.Range("F6").Value  = 419
$ws4.Range("F9").Value  = 15460
$ws4.Range("F11").Value = 9065
$ws4.Range("F12").Value = 389
$ws4.Range("F15").Value = 102
$ws4.Range("F18").Value = 204
$ws4.Range("F32").Value = 70
$ws4.Range("F37").Value = 259
